# Append a new log entry row to the end of the progress table:
#   Date: 24/03/2025
#   What I achieved: Improving jumping physics, making them feel more natural

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Add a new row at the bottom of the table; it inherits the formatting
# (column widths, row height, centered date column) from the preceding row.
$newRow = $table.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "24/03/2025"
$newRow.Cells.Item(2).Range.Text = "Improving jumping physics, making them feel more natural"
